# Fixed bug in models
#
# 1. "pars" sheet: update default values for the incubation (i) and
#    recovery (r) rate constants.
$wb = $excel.ActiveWorkbook

$parsWs = $wb.Worksheets.Item("pars")
$parsWs.Range("D3").Value = 0.2
$parsWs.Range("D4").Value = 0.4

# 2. "eqns" sheet: the stoichiometry columns (S/E/I/R, columns E:H) held
#    text placeholders such as "-t"/"t" instead of actual numeric
#    coefficients. Replace them with the real numbers.
$eqnsWs = $wb.Worksheets.Item("eqns")
$eqnsWs.Range("E2").Value = -1
$eqnsWs.Range("F2").Value = 1

$eqnsWs.Range("F3").Value = -1
$eqnsWs.Range("G3").Value = 1

$eqnsWs.Range("G4").Value = -1
$eqnsWs.Range("H4").Value = 1

# 3. Move the active selection/tab: "pars" becomes the active sheet
#    (previously "eqns" was active).
[void]$eqnsWs.Range("F18").Select()

[void]$parsWs.Activate()
[void]$parsWs.Range("D5").Select()
